$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3362.25
$ws.Range("I43").Value = 3527.2727
$ws.Range("K43").Value = 3527.2727
$ws.Range("M43").Value = -3458.2727
$ws.Range("H98").Value = 1295.6097
$ws.Range("I98").Value = 995.5405
$ws.Range("J98").Value = 4071.25
$ws.Range("K98").Value = 995.5405
$ws.Range("L98").Value = 4071.25
$ws.Range("M98").Value = 502.4595
$ws.Range("N98").Value = -7067.25
$ws.Range("H103").Value = 537.5
$ws.Range("I103").Value = 450
$ws.Range("J103").Value = 625
$ws.Range("K103").Value = 1350
$ws.Range("L103").Value = 1875
$ws.Range("M103").Value = -764
$ws.Range("N103").Value = -3047
$ws.Range("H122").Value = 1295.6097
$ws.Range("I122").Value = 995.5405
$ws.Range("J122").Value = 4071.25
$ws.Range("K122").Value = 2986.6215
$ws.Range("L122").Value = 12213.75
$ws.Range("M122").Value = -536.6214999999997
$ws.Range("N122").Value = -17113.75
$ws.Range("H137").Value = 4992.4346
$ws.Range("I137").Value = 2398.5881
$ws.Range("J137").Value = 12341.667
$ws.Range("K137").Value = 7195.7643
$ws.Range("L137").Value = 37025.001
$ws.Range("M137").Value = -4645.7643
$ws.Range("N137").Value = -42125.001
$ws.Range("H138").Value = 6466.449
$ws.Range("J138").Value = 6194.7075
$ws.Range("L138").Value = 18584.1225
$ws.Range("N138").Value = -28864.1225

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 48968.43
$ws.Range("I2").Value = 60041.234
$ws.Range("K2").Value = 60041.234
$ws.Range("M2").Value = -59928.234
$ws.Range("H5").Value = 74.15000000000001
$ws.Range("J5").Value = 74.2
$ws.Range("L5").Value = 74.2
$ws.Range("N5").Value = -298.2
$ws.Range("H32").Value = 3537
$ws.Range("I32").Value = 2532.691
$ws.Range("K32").Value = 2532.691
$ws.Range("M32").Value = -2245.691
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H116").Value = 48968.43
$ws.Range("I116").Value = 60041.234
$ws.Range("K116").Value = 60041.234
$ws.Range("M116").Value = -57747.234

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 48968.43
$ws.Range("I3").Value = 60041.234
$ws.Range("K3").Value = 60041.234
$ws.Range("M3").Value = -59927.234
$ws.Range("H4").Value = 74.15000000000001
$ws.Range("J4").Value = 74.2
$ws.Range("L4").Value = 74.2
$ws.Range("N4").Value = -304.2
$ws.Range("H20").Value = 2489.5454
$ws.Range("I20").Value = 1581.1666
$ws.Range("K20").Value = 1581.1666
$ws.Range("M20").Value = -1334.1666
$ws.Range("H80").Value = 307.45
$ws.Range("I80").Value = 42.857143
$ws.Range("J80").Value = 449.92307
$ws.Range("K80").Value = 42.857143
$ws.Range("L80").Value = 449.92307
$ws.Range("M80").Value = 955.142857
$ws.Range("N80").Value = -2445.92307
$ws.Range("H82").Value = 33639.223
$ws.Range("J82").Value = 105300
$ws.Range("L82").Value = 105300
$ws.Range("N82").Value = -106066
$ws.Range("H83").Value = 307.45
$ws.Range("I83").Value = 42.857143
$ws.Range("J83").Value = 449.92307
$ws.Range("K83").Value = 214.285715
$ws.Range("L83").Value = 2249.61535
$ws.Range("M83").Value = 4777.714285
$ws.Range("N83").Value = -12233.61535
$ws.Range("H85").Value = 33639.223
$ws.Range("J85").Value = 105300
$ws.Range("L85").Value = 105300
$ws.Range("N85").Value = -107952
$ws.Range("H99").Value = 3692.5334
$ws.Range("I99").Value = 3427.5454
$ws.Range("J99").Value = 4421.25
$ws.Range("K99").Value = 3427.5454
$ws.Range("L99").Value = 4421.25
$ws.Range("M99").Value = -1929.5454
$ws.Range("N99").Value = -7417.25
$ws.Range("H134").Value = 51782.227
$ws.Range("I134").Value = 4974
$ws.Range("J134").Value = 98590.45
$ws.Range("K134").Value = 14922
$ws.Range("L134").Value = 295771.35
$ws.Range("M134").Value = -12387
$ws.Range("N134").Value = -300841.35

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 594.9524
$ws.Range("I7").Value = 526.3570999999999
$ws.Range("J7").Value = 732.1429000000001
$ws.Range("K7").Value = 526.3570999999999
$ws.Range("L7").Value = 732.1429000000001
$ws.Range("M7").Value = -413.3570999999999
$ws.Range("N7").Value = -958.1429000000001
$ws.Range("H16").Value = 4089.244
$ws.Range("J16").Value = 8874
$ws.Range("L16").Value = 8874
$ws.Range("N16").Value = -9448
$ws.Range("H58").Value = 203607.3
$ws.Range("I58").Value = 372768.22
$ws.Range("J58").Value = 5027.087
$ws.Range("K58").Value = 372768.22
$ws.Range("L58").Value = 5027.087
$ws.Range("M58").Value = -372565.22
$ws.Range("N58").Value = -5433.087
$ws.Range("H99").Value = 8731.157999999999
$ws.Range("J99").Value = 7087.5
$ws.Range("L99").Value = 7087.5
$ws.Range("N99").Value = -10083.5
$ws.Range("H103").Value = 25999.4
$ws.Range("I103").Value = 27499.75
$ws.Range("K103").Value = 27499.75
$ws.Range("M103").Value = -26327.75
$ws.Range("H113").Value = 4089.244
$ws.Range("J113").Value = 8874
$ws.Range("L113").Value = 8874
$ws.Range("N113").Value = -13214
$ws.Range("H126").Value = 8731.157999999999
$ws.Range("J126").Value = 7087.5
$ws.Range("L126").Value = 21262.5
$ws.Range("N126").Value = -26202.5
$ws.Range("H134").Value = 212267.31
$ws.Range("I134").Value = 2856.1155
$ws.Range("K134").Value = 8568.3465
$ws.Range("M134").Value = -6033.3465
$ws.Range("H136").Value = 203607.3
$ws.Range("I136").Value = 372768.22
$ws.Range("J136").Value = 5027.087
$ws.Range("K136").Value = 1118304.66
$ws.Range("L136").Value = 15081.261
$ws.Range("M136").Value = -1115754.66
$ws.Range("N136").Value = -20181.261

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 3000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8540
$ws.Range("H55").Value = 9859
$ws.Range("I55").Value = 1264.1428
$ws.Range("K55").Value = 3792.4284
$ws.Range("M55").Value = -3615.4284
$ws.Range("H75").Value = 1900
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 1900
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H98").Value = 2028.5714
$ws.Range("I98").Value = 3554.5
$ws.Range("J98").Value = 1418.2
$ws.Range("K98").Value = 10663.5
$ws.Range("L98").Value = 4254.6
$ws.Range("M98").Value = -9165.5
$ws.Range("N98").Value = -7250.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3580.7097
$ws.Range("I46").Value = 2605.2632
$ws.Range("K46").Value = 2605.2632
$ws.Range("M46").Value = -2417.2632
$ws.Range("H132").Value = 3585.158
$ws.Range("I132").Value = 1215.5714
$ws.Range("J132").Value = 6512.294
$ws.Range("K132").Value = 3646.7142
$ws.Range("L132").Value = 19536.882
$ws.Range("M132").Value = -1116.7142
$ws.Range("N132").Value = -24596.882

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 84292.35000000001
$ws.Range("I40").Value = 84295.92
$ws.Range("J40").Value = 84285.71000000001
$ws.Range("K40").Value = 84295.92
$ws.Range("L40").Value = 84285.71000000001
$ws.Range("M40").Value = -84146.92
$ws.Range("N40").Value = -84583.71000000001
$ws.Range("H107").Value = 322.3684
$ws.Range("J107").Value = 444
$ws.Range("L107").Value = 1332
$ws.Range("N107").Value = -5172
$ws.Range("H136").Value = 67031.98
$ws.Range("I136").Value = 14631.871
$ws.Range("K136").Value = 43895.613
$ws.Range("M136").Value = -41345.613
